# Auto-generated edit script: updates cached market-price / profit values
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# the latest Universalis marketboard snapshot pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: columns H,I,J,K,L,M,N
$ws.Range("H28").Value = 520.70966
$ws.Range("I28").Value = 585
$ws.Range("J28").Value = 300.2857
$ws.Range("K28").Value = 585
$ws.Range("L28").Value = 300.2857
$ws.Range("M28").Value = -100
$ws.Range("N28").Value = -1270.2857
# Row 69: columns H,I,J,K,L,M,N
$ws.Range("H69").Value = 3839.25
$ws.Range("I69").Value = 4000
$ws.Range("J69").Value = 3785.6667
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 11357.0001
$ws.Range("M69").Value = -11126
$ws.Range("N69").Value = -13105.0001
# Row 72: columns H,I,J,K,L,M,N
$ws.Range("H72").Value = 3839.25
$ws.Range("I72").Value = 4000
$ws.Range("J72").Value = 3785.6667
$ws.Range("K72").Value = 36000
$ws.Range("L72").Value = 34071.0003
$ws.Range("M72").Value = -31632
$ws.Range("N72").Value = -42807.0003
# Row 107: columns H,I,J,K,L,M,N
$ws.Range("H107").Value = 2604681.8
$ws.Range("I107").Value = 3788469.8
$ws.Range("J107").Value = 348.4
$ws.Range("K107").Value = 3788469.8
$ws.Range("L107").Value = 348.4
$ws.Range("M107").Value = -3786549.8
$ws.Range("N107").Value = -4188.4
# Row 116: columns H,I,J,K,L,M,N
$ws.Range("H116").Value = 4574.303
$ws.Range("I116").Value = 5117.3335
$ws.Range("J116").Value = 4121.778
$ws.Range("K116").Value = 5117.3335
$ws.Range("L116").Value = 4121.778
$ws.Range("M116").Value = -1675.3335
$ws.Range("N116").Value = -11005.778
# Row 132: columns H,I,K,M
$ws.Range("H132").Value = 1762.2084
$ws.Range("I132").Value = 1251.1765
$ws.Range("K132").Value = 3753.5295
$ws.Range("M132").Value = -1223.5295
# Row 137: columns H,I,J,K,L,M,N
$ws.Range("H137").Value = 3948634.2
$ws.Range("I137").Value = 1725376
$ws.Range("J137").Value = 11112467
$ws.Range("K137").Value = 5176128
$ws.Range("L137").Value = 33337401
$ws.Range("M137").Value = -5173578
$ws.Range("N137").Value = -33342501

$ws = $wb.Worksheets.Item("ARM")
# Row 32: columns H,I,J,K,L,M,N
$ws.Range("H32").Value = 16882.58
$ws.Range("I32").Value = 4111.136
$ws.Range("J32").Value = 110539.836
$ws.Range("K32").Value = 4111.136
$ws.Range("L32").Value = 110539.836
$ws.Range("M32").Value = -3824.136
$ws.Range("N32").Value = -111113.836
# Row 74: columns H,I,J,K,L,M,N
$ws.Range("H74").Value = 903.193
$ws.Range("I74").Value = 866.9808
$ws.Range("J74").Value = 1279.8
$ws.Range("K74").Value = 866.9808
$ws.Range("L74").Value = 1279.8
$ws.Range("M74").Value = 7.019199999999955
$ws.Range("N74").Value = -3027.8
# Row 77: columns H,I,J,K,L,M,N
$ws.Range("H77").Value = 903.193
$ws.Range("I77").Value = 866.9808
$ws.Range("J77").Value = 1279.8
$ws.Range("K77").Value = 4334.904
$ws.Range("L77").Value = 6399
$ws.Range("M77").Value = 33.09599999999955
$ws.Range("N77").Value = -15135
# Row 122: columns H,I,J,K,L,M,N
$ws.Range("H122").Value = 2136.3684
$ws.Range("I122").Value = 2042.2142
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 6126.642599999999
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -3676.642599999999
$ws.Range("N122").Value = -12100

$ws = $wb.Worksheets.Item("BSM")
# Row 94: columns H,I,J,K,L,M,N
$ws.Range("H94").Value = 446.19354
$ws.Range("I94").Value = 440.36365
$ws.Range("J94").Value = 460.44446
$ws.Range("K94").Value = 440.36365
$ws.Range("L94").Value = 460.44446
$ws.Range("M94").Value = 10.63634999999999
$ws.Range("N94").Value = -1362.44446

$ws = $wb.Worksheets.Item("CRP")
# Row 31: columns H,I,J,K,L,M,N
$ws.Range("H31").Value = 2258.1724
$ws.Range("I31").Value = 1818.6364
$ws.Range("J31").Value = 3639.5715
$ws.Range("K31").Value = 1818.6364
$ws.Range("L31").Value = 3639.5715
$ws.Range("M31").Value = -1523.6364
$ws.Range("N31").Value = -4229.5715
# Row 34: columns H,I,J,K,L,M,N
$ws.Range("H34").Value = 2258.1724
$ws.Range("I34").Value = 1818.6364
$ws.Range("J34").Value = 3639.5715
$ws.Range("K34").Value = 1818.6364
$ws.Range("L34").Value = 3639.5715
$ws.Range("M34").Value = -1616.6364
$ws.Range("N34").Value = -4043.5715
# Row 132: columns H,I,J,K,L,M,N
$ws.Range("H132").Value = 4088.6924
$ws.Range("I132").Value = 3525.25
$ws.Range("J132").Value = 5966.8335
$ws.Range("K132").Value = 10575.75
$ws.Range("L132").Value = 17900.5005
$ws.Range("M132").Value = -8045.75
$ws.Range("N132").Value = -22960.5005
# Row 138: columns H,J,L,N
$ws.Range("H138").Value = 46820
$ws.Range("J138").Value = 46820
$ws.Range("L138").Value = 46820
$ws.Range("N138").Value = -57100

$ws = $wb.Worksheets.Item("CUL")
# Row 80: columns H,I,J,K,L,M,N
$ws.Range("H80").Value = 6108.077
$ws.Range("I80").Value = 2702
$ws.Range("J80").Value = 6391.9165
$ws.Range("K80").Value = 8106
$ws.Range("L80").Value = 19175.7495
$ws.Range("M80").Value = -7170
$ws.Range("N80").Value = -21047.7495
# Row 83: columns H,I,J,K,L,M,N
$ws.Range("H83").Value = 6108.077
$ws.Range("I83").Value = 2702
$ws.Range("J83").Value = 6391.9165
$ws.Range("K83").Value = 24318
$ws.Range("L83").Value = 57527.2485
$ws.Range("M83").Value = -19638
$ws.Range("N83").Value = -66887.2485
# Row 139: columns H,J,L,N
$ws.Range("H139").Value = 25164.93
$ws.Range("J139").Value = 170005.33
$ws.Range("L139").Value = 510015.99
$ws.Range("N139").Value = -520295.99
# Row 140: columns H,I,K,M
$ws.Range("H140").Value = 31085.777
$ws.Range("I140").Value = 67741.336
$ws.Range("K140").Value = 203224.008
$ws.Range("M140").Value = -198044.008

$ws = $wb.Worksheets.Item("GSM")
# Row 15: columns H,J,L,N
$ws.Range("H15").Value = 39000
$ws.Range("J15").Value = 39000
$ws.Range("L15").Value = 39000
$ws.Range("N15").Value = -39576
# Row 81: columns H,J,L,N
$ws.Range("H81").Value = 39000
$ws.Range("J81").Value = 39000
$ws.Range("L81").Value = 39000
$ws.Range("N81").Value = -40996
# Row 84: columns H,J,L,N
$ws.Range("H84").Value = 39000
$ws.Range("J84").Value = 39000
$ws.Range("L84").Value = 117000
$ws.Range("N84").Value = -126984
# Row 122: columns H,I,J,K,L,M,N
$ws.Range("H122").Value = 2515.3572
$ws.Range("I122").Value = 2145.2222
$ws.Range("J122").Value = 3181.6
$ws.Range("K122").Value = 6435.6666
$ws.Range("L122").Value = 9544.799999999999
$ws.Range("M122").Value = -3985.6666
$ws.Range("N122").Value = -14444.8
# Row 123: columns H,J,L,N
$ws.Range("H123").Value = 13876.583
$ws.Range("J123").Value = 13876.583
$ws.Range("L123").Value = 13876.583
$ws.Range("N123").Value = -18776.583

$ws = $wb.Worksheets.Item("LTW")
# Row 7: columns H,I,J,K,L,M,N
$ws.Range("H7").Value = 5265544
$ws.Range("I7").Value = 9092946
$ws.Range("J7").Value = 2866.875
$ws.Range("K7").Value = 9092946
$ws.Range("L7").Value = 2866.875
$ws.Range("M7").Value = -9092834
$ws.Range("N7").Value = -3090.875
# Row 61: columns H,I,K,M
$ws.Range("H61").Value = 1823.05
$ws.Range("I61").Value = 1597.4117
$ws.Range("K61").Value = 1597.4117
$ws.Range("M61").Value = -1395.4117
# Row 113: columns H,I,K,M
$ws.Range("H113").Value = 1823.05
$ws.Range("I113").Value = 1597.4117
$ws.Range("K113").Value = 1597.4117
$ws.Range("M113").Value = 572.5882999999999
# Row 126: columns H,I,J,K,L,M,N
$ws.Range("H126").Value = 5265544
$ws.Range("I126").Value = 9092946
$ws.Range("J126").Value = 2866.875
$ws.Range("K126").Value = 27278838
$ws.Range("L126").Value = 8600.625
$ws.Range("M126").Value = -27276368
$ws.Range("N126").Value = -13540.625
# Row 132: columns H,I,J,K,L,M,N
$ws.Range("H132").Value = 1902.8572
$ws.Range("I132").Value = 1732.7273
$ws.Range("J132").Value = 3400
$ws.Range("K132").Value = 5198.1819
$ws.Range("L132").Value = 10200
$ws.Range("M132").Value = -2668.1819
$ws.Range("N132").Value = -15260

$ws = $wb.Worksheets.Item("WVR")
# Row 75: columns H,J,L,N
$ws.Range("H75").Value = 24658.182
$ws.Range("J75").Value = 24658.182
$ws.Range("L75").Value = 24658.182
$ws.Range("N75").Value = -26530.182
# Row 78: columns H,J,L,N
$ws.Range("H78").Value = 24658.182
$ws.Range("J78").Value = 24658.182
$ws.Range("L78").Value = 73974.546
$ws.Range("N78").Value = -83334.546
# Row 132: columns H,I,K,M
$ws.Range("H132").Value = 1472.5312
$ws.Range("I132").Value = 2163.5
$ws.Range("K132").Value = 6490.5
$ws.Range("M132").Value = -3960.5
# Row 136: columns H,I,J,K,L,M,N
$ws.Range("H136").Value = 1614.0546
$ws.Range("I136").Value = 1677.9556
$ws.Range("J136").Value = 1326.5
$ws.Range("K136").Value = 5033.8668
$ws.Range("L136").Value = 3979.5
$ws.Range("M136").Value = -2483.8668
$ws.Range("N136").Value = -9079.5

Write-Host "Updated 199 cells across 8 sheets"